# Atualizacao de bases das ligas, do dia: 11-04-2024 as 00:31
#
# The re-pulled feed re-ordered the two fixtures played on each of three
# match days, so the full data row (every column except the running "id"
# index in column A) needs to be swapped between:
#   - rows 2   and 3   (2023-07-29 fixtures)
#   - rows 167 and 168 (2024-02-18 fixtures)
#   - rows 195 and 196 (2024-03-12 fixtures)
#
# Column A (the sequential 0-based id) stays put on each row; columns
# B:AC (id/match-code, teams, scores, odds, ...) swap as a block.
#
# NOTE: this host's `.Value` getter returns a reflection stub instead of
# the live cell content, so `.Value2` is used to round-trip the row data
# between ranges (the `.Value` setter itself works fine and is used for
# plain scalar writes elsewhere).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowPairs = @(@(2, 3), @(167, 168), @(195, 196))

foreach ($pair in $rowPairs) {
    $rowA = $pair[0]
    $rowB = $pair[1]

    $rangeA = $ws.Range("B$rowA" + ":AC$rowA")
    $rangeB = $ws.Range("B$rowB" + ":AC$rowB")

    $dataA = $rangeA.Value2
    $dataB = $rangeB.Value2

    $rangeA.Value2 = $dataB
    $rangeB.Value2 = $dataA
}
